# Add a "Comments" column (R) to Sheet1, marking archived job rows with "Archive".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cell: bold font + yellow fill (matches the other bold header cells' font,
# plus a new highlight fill).
$ws.Range("R1").Value = "Comments"
$ws.Range("R1").Font.Bold = $true
$ws.Range("R1").Interior.Color = 65535

# Data rows flagged as archived jobs.
$archiveRanges = @(
    "R2:R3",
    "R11:R36",
    "R40:R42",
    "R55:R68",
    "R70:R95",
    "R97:R103",
    "R111:R118"
)

foreach ($rng in $archiveRanges) {
    $ws.Range($rng).Value = "Archive"
}

# Reset the view: zoom out and move the selection back to the top of the sheet
# (matches the saved sheetView state after the edit).
$excel.ActiveWindow.Zoom = 80
[void]$ws.Range("B2").Select()
